$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# Sheet 2: Главные
$ws2.Range("C2").Value = 21
$ws2.Range("D2").Value = 498
$ws2.Range("E2").Value = 216
$ws2.Range("F2").Value = 282
$ws2.Range("G2").Value = 23.71
$ws2.Range("H2").Value = 10.29
$ws2.Range("I2").Value = 13.43
$ws2.Range("J2").Value = 93
$ws2.Range("K2").Value = 111
$ws2.Range("W2").Value = 14
$ws2.Range("AA2").Value = "2025-11-06 03:04:31"
$ws2.Range("C3").Value = 21
$ws2.Range("D3").Value = 381
$ws2.Range("E3").Value = 167
$ws2.Range("F3").Value = 214
$ws2.Range("G3").Value = 18.14
$ws2.Range("I3").Value = 10.19
$ws2.Range("J3").Value = 81
$ws2.Range("K3").Value = 87
$ws2.Range("AA3").Value = "2025-11-06 03:04:31"
$ws2.Range("C4").Value = 16
$ws2.Range("D4").Value = 278
$ws2.Range("E4").Value = 121
$ws2.Range("F4").Value = 157
$ws2.Range("G4").Value = 17.38
$ws2.Range("H4").Value = 7.56
$ws2.Range("I4").Value = 9.81
$ws2.Range("J4").Value = 58
$ws2.Range("K4").Value = 66
$ws2.Range("W4").Value = 10
$ws2.Range("AA4").Value = "2025-11-06 03:04:31"
$ws2.Range("C5").Value = 22
$ws2.Range("D5").Value = 349
$ws2.Range("E5").Value = 191
$ws2.Range("F5").Value = 158
$ws2.Range("G5").Value = 15.86
$ws2.Range("H5").Value = 8.68
$ws2.Range("I5").Value = 7.18
$ws2.Range("J5").Value = 93
$ws2.Range("K5").Value = 79
$ws2.Range("V5").Value = 16
$ws2.Range("AA5").Value = "2025-11-06 03:04:31"
$ws2.Range("C6").Value = 21
$ws2.Range("D6").Value = 393
$ws2.Range("E6").Value = 165
$ws2.Range("F6").Value = 228
$ws2.Range("G6").Value = 18.71
$ws2.Range("H6").Value = 7.86
$ws2.Range("I6").Value = 10.86
$ws2.Range("J6").Value = 75
$ws2.Range("K6").Value = 94
$ws2.Range("AA6").Value = "2025-11-06 03:04:31"
$ws2.Range("AA7").Value = "2025-11-06 03:04:31"
$ws2.Range("C8").Value = 20
$ws2.Range("D8").Value = 378
$ws2.Range("E8").Value = 185
$ws2.Range("F8").Value = 193
$ws2.Range("G8").Value = 18.9
$ws2.Range("H8").Value = 9.25
$ws2.Range("I8").Value = 9.65
$ws2.Range("J8").Value = 85
$ws2.Range("K8").Value = 89
$ws2.Range("AA8").Value = "2025-11-06 03:04:31"
$ws2.Range("C9").Value = 22
$ws2.Range("D9").Value = 356
$ws2.Range("E9").Value = 193
$ws2.Range("F9").Value = 163
$ws2.Range("G9").Value = 16.18
$ws2.Range("H9").Value = 8.77
$ws2.Range("I9").Value = 7.41
$ws2.Range("J9").Value = 94
$ws2.Range("K9").Value = 79
$ws2.Range("AA9").Value = "2025-11-06 03:04:31"
$ws2.Range("AA10").Value = "2025-11-06 03:04:31"
$ws2.Range("AA11").Value = "2025-11-06 03:04:31"
$ws2.Range("AA12").Value = "2025-11-06 03:04:31"
$ws2.Range("C13").Value = 8
$ws2.Range("D13").Value = 135
$ws2.Range("E13").Value = 72
$ws2.Range("F13").Value = 63
$ws2.Range("G13").Value = 16.88
$ws2.Range("H13").Value = 9
$ws2.Range("I13").Value = 7.88
$ws2.Range("J13").Value = 36
$ws2.Range("K13").Value = 29
$ws2.Range("T13").Value = 1
$ws2.Range("AA13").Value = "2025-11-06 03:04:31"
$ws2.Range("AA14").Value = "2025-11-06 03:04:31"
$ws2.Range("C15").Value = 14
$ws2.Range("D15").Value = 200
$ws2.Range("E15").Value = 90
$ws2.Range("F15").Value = 110
$ws2.Range("G15").Value = 14.29
$ws2.Range("H15").Value = 6.43
$ws2.Range("I15").Value = 7.86
$ws2.Range("J15").Value = 45
$ws2.Range("K15").Value = 55
$ws2.Range("AA15").Value = "2025-11-06 03:04:31"
$ws2.Range("C16").Value = 21
$ws2.Range("D16").Value = 430
$ws2.Range("E16").Value = 207
$ws2.Range("F16").Value = 223
$ws2.Range("G16").Value = 20.48
$ws2.Range("H16").Value = 9.86
$ws2.Range("I16").Value = 10.62
$ws2.Range("J16").Value = 76
$ws2.Range("K16").Value = 79
$ws2.Range("AA16").Value = "2025-11-06 03:04:31"
$ws2.Range("AA17").Value = "2025-11-06 03:04:31"
$ws2.Range("AA18").Value = "2025-11-06 03:04:31"
$ws2.Range("AA19").Value = "2025-11-06 03:04:31"
$ws2.Range("AA20").Value = "2025-11-06 03:04:31"
$ws2.Range("C21").Value = 17
$ws2.Range("D21").Value = 259
$ws2.Range("E21").Value = 110
$ws2.Range("F21").Value = 149
$ws2.Range("G21").Value = 15.24
$ws2.Range("H21").Value = 6.47
$ws2.Range("I21").Value = 8.76
$ws2.Range("J21").Value = 45
$ws2.Range("K21").Value = 62
$ws2.Range("AA21").Value = "2025-11-06 03:04:31"
$ws2.Range("C22").Value = 15
$ws2.Range("D22").Value = 294
$ws2.Range("E22").Value = 116
$ws2.Range("F22").Value = 178
$ws2.Range("G22").Value = 19.6
$ws2.Range("H22").Value = 7.73
$ws2.Range("I22").Value = 11.87
$ws2.Range("J22").Value = 58
$ws2.Range("K22").Value = 59
$ws2.Range("V22").Value = 4
$ws2.Range("AA22").Value = "2025-11-06 03:04:31"
$ws2.Range("AA23").Value = "2025-11-06 03:04:31"
$ws2.Range("AA24").Value = "2025-11-06 03:04:31"
$ws2.Range("C25").Value = 21
$ws2.Range("D25").Value = 372
$ws2.Range("E25").Value = 190
$ws2.Range("F25").Value = 182
$ws2.Range("G25").Value = 17.71
$ws2.Range("H25").Value = 9.05
$ws2.Range("I25").Value = 8.67
$ws2.Range("J25").Value = 90
$ws2.Range("K25").Value = 86
$ws2.Range("AA25").Value = "2025-11-06 03:04:31"
$ws2.Range("AA26").Value = "2025-11-06 03:04:31"

# Sheet 3: Линейные
$ws3.Range("AA2").Value = "2025-11-06 03:04:31"
$ws3.Range("AA3").Value = "2025-11-06 03:04:31"
$ws3.Range("C4").Value = 11
$ws3.Range("D4").Value = 182
$ws3.Range("E4").Value = 78
$ws3.Range("F4").Value = 104
$ws3.Range("G4").Value = 16.55
$ws3.Range("H4").Value = 7.09
$ws3.Range("I4").Value = 9.45
$ws3.Range("J4").Value = 39
$ws3.Range("K4").Value = 42
$ws3.Range("V4").Value = 6
$ws3.Range("AA4").Value = "2025-11-06 03:04:31"
$ws3.Range("C5").Value = 9
$ws3.Range("D5").Value = 132
$ws3.Range("E5").Value = 68
$ws3.Range("F5").Value = 64
$ws3.Range("G5").Value = 14.67
$ws3.Range("H5").Value = 7.56
$ws3.Range("I5").Value = 7.11
$ws3.Range("J5").Value = 34
$ws3.Range("K5").Value = 32
$ws3.Range("T5").Value = 1
$ws3.Range("AA5").Value = "2025-11-06 03:04:31"
$ws3.Range("AA6").Value = "2025-11-06 03:04:31"
$ws3.Range("C7").Value = 12
$ws3.Range("D7").Value = 209
$ws3.Range("E7").Value = 72
$ws3.Range("F7").Value = 137
$ws3.Range("G7").Value = 17.42
$ws3.Range("H7").Value = 6
$ws3.Range("I7").Value = 11.42
$ws3.Range("J7").Value = 36
$ws3.Range("K7").Value = 41
$ws3.Range("AA7").Value = "2025-11-06 03:04:31"
$ws3.Range("AA8").Value = "2025-11-06 03:04:31"
$ws3.Range("AA9").Value = "2025-11-06 03:04:31"
$ws3.Range("AA10").Value = "2025-11-06 03:04:31"
$ws3.Range("C11").Value = 12
$ws3.Range("D11").Value = 171
$ws3.Range("E11").Value = 84
$ws3.Range("F11").Value = 87
$ws3.Range("G11").Value = 14.25
$ws3.Range("H11").Value = 7
$ws3.Range("I11").Value = 7.25
$ws3.Range("J11").Value = 42
$ws3.Range("K11").Value = 41
$ws3.Range("T11").Value = 1
$ws3.Range("AA11").Value = "2025-11-06 03:04:31"
$ws3.Range("AA12").Value = "2025-11-06 03:04:31"
$ws3.Range("AA13").Value = "2025-11-06 03:04:31"
$ws3.Range("AA14").Value = "2025-11-06 03:04:31"
$ws3.Range("C15").Value = 18
$ws3.Range("D15").Value = 375
$ws3.Range("E15").Value = 195
$ws3.Range("F15").Value = 180
$ws3.Range("G15").Value = 20.83
$ws3.Range("H15").Value = 10.83
$ws3.Range("I15").Value = 10
$ws3.Range("J15").Value = 75
$ws3.Range("K15").Value = 70
$ws3.Range("AA15").Value = "2025-11-06 03:04:31"
$ws3.Range("AA16").Value = "2025-11-06 03:04:31"
$ws3.Range("AA17").Value = "2025-11-06 03:04:31"
$ws3.Range("C18").Value = 22
$ws3.Range("D18").Value = 391
$ws3.Range("E18").Value = 188
$ws3.Range("F18").Value = 203
$ws3.Range("G18").Value = 17.77
$ws3.Range("H18").Value = 8.55
$ws3.Range("I18").Value = 9.23
$ws3.Range("J18").Value = 89
$ws3.Range("K18").Value = 84
$ws3.Range("AA18").Value = "2025-11-06 03:04:31"
$ws3.Range("C19").Value = 19
$ws3.Range("D19").Value = 349
$ws3.Range("E19").Value = 162
$ws3.Range("F19").Value = 187
$ws3.Range("G19").Value = 18.37
$ws3.Range("H19").Value = 8.53
$ws3.Range("I19").Value = 9.84
$ws3.Range("J19").Value = 76
$ws3.Range("K19").Value = 81
$ws3.Range("AA19").Value = "2025-11-06 03:04:31"
$ws3.Range("AA20").Value = "2025-11-06 03:04:31"
$ws3.Range("AA21").Value = "2025-11-06 03:04:31"
$ws3.Range("C22").Value = 16
$ws3.Range("D22").Value = 259
$ws3.Range("E22").Value = 124
$ws3.Range("F22").Value = 135
$ws3.Range("G22").Value = 16.19
$ws3.Range("H22").Value = 7.75
$ws3.Range("I22").Value = 8.44
$ws3.Range("J22").Value = 62
$ws3.Range("K22").Value = 65
$ws3.Range("W22").Value = 20
$ws3.Range("AA22").Value = "2025-11-06 03:04:31"
$ws3.Range("AA23").Value = "2025-11-06 03:04:31"
$ws3.Range("C24").Value = 22
$ws3.Range("D24").Value = 404
$ws3.Range("E24").Value = 163
$ws3.Range("F24").Value = 241
$ws3.Range("G24").Value = 18.36
$ws3.Range("H24").Value = 7.41
$ws3.Range("I24").Value = 10.95
$ws3.Range("J24").Value = 74
$ws3.Range("K24").Value = 93
$ws3.Range("AA24").Value = "2025-11-06 03:04:31"
$ws3.Range("AA25").Value = "2025-11-06 03:04:31"
$ws3.Range("C26").Value = 19
$ws3.Range("D26").Value = 414
$ws3.Range("E26").Value = 175
$ws3.Range("F26").Value = 239
$ws3.Range("G26").Value = 21.79
$ws3.Range("H26").Value = 9.21
$ws3.Range("I26").Value = 12.58
$ws3.Range("J26").Value = 65
$ws3.Range("K26").Value = 67
$ws3.Range("AA26").Value = "2025-11-06 03:04:31"
